$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.3424902588620627
$ws.Range("J2").Value = 0.3424902588620627
$ws.Range("M2").Value = 0.1036536666666667
$ws.Range("N2").Value = 0.310961
$ws.Range("O2").Value = 0.08145093039891602
$ws.Range("P2").Value = 0.08145093039891602
$ws.Range("Q2").Value = 0.1395158313624444
$ws.Range("R2").Value = 1.255642482262
$ws.Range("S2").Value = 0.0278961502368806
$ws.Range("T2").Value = 0.0278961502368806

# Row 3
$ws.Range("I3").Value = 0.3424902588620627
$ws.Range("J3").Value = 0.3424902588620627
$ws.Range("O3").Value = 0.7918149097994615
$ws.Range("P3").Value = 0.7918149097994616
$ws.Range("S3").Value = 0.2711888934280584
$ws.Range("T3").Value = 0.2711888934280585

# Row 4
$ws.Range("I4").Value = 0.3424902588620627
$ws.Range("J4").Value = 0.3424902588620627
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1612806666666667
$ws.Range("N4").Value = 0.483842
$ws.Range("O4").Value = 0.1267341598016225
$ws.Range("P4").Value = 0.1267341598016225
$ws.Range("Q4").Value = 0.2170806592404444
$ws.Range("R4").Value = 1.953725933164
$ws.Range("S4").Value = 0.04340521519712371
$ws.Range("T4").Value = 0.04340521519712371

# Row 5
$ws.Range("G5").Value = 2.584001666666667
$ws.Range("H5").Value = 7.752005
$ws.Range("I5").Value = 0.6575097411379373
$ws.Range("J5").Value = 0.6575097411379373
$ws.Range("M5").Value = 0.1036536666666667
$ws.Range("N5").Value = 0.310961
$ws.Range("O5").Value = 0.08145093039891602
$ws.Range("P5").Value = 0.08145093039891602
$ws.Range("Q5").Value = 0.2678412474227778
$ws.Range("R5").Value = 2.410571226805
$ws.Range("S5").Value = 0.05355478016203542
$ws.Range("T5").Value = 0.05355478016203542

# Row 6
$ws.Range("G6").Value = 2.584001666666667
$ws.Range("H6").Value = 7.752005
$ws.Range("I6").Value = 0.6575097411379373
$ws.Range("J6").Value = 0.6575097411379373
$ws.Range("O6").Value = 0.7918149097994615
$ws.Range("P6").Value = 0.7918149097994616
$ws.Range("Q6").Value = 2.603784783426666
$ws.Range("R6").Value = 23.43406305084
$ws.Range("S6").Value = 0.520626016371403
$ws.Range("T6").Value = 0.5206260163714032

# Row 7
$ws.Range("G7").Value = 2.584001666666667
$ws.Range("H7").Value = 7.752005
$ws.Range("I7").Value = 0.6575097411379373
$ws.Range("J7").Value = 0.6575097411379373
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1612806666666667
$ws.Range("N7").Value = 0.483842
$ws.Range("O7").Value = 0.1267341598016225
$ws.Range("P7").Value = 0.1267341598016225
$ws.Range("Q7").Value = 0.4167495114677777
$ws.Range("R7").Value = 3.75074560321
$ws.Range("S7").Value = 0.08332894460449877
$ws.Range("T7").Value = 0.08332894460449877
